# Update bulk list upload template:
# - Rename the "Primer Nombre" column header to just "Nombre"
# - Leave the active selection on B5 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Nombre"

$ws.Range("B5").Select()
